$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells hold plain text such as "60.822.60" or "41.40".
# If we just assign a numeric-looking string, Excel/COM will happily turn it
# into a floating point number and mangle/round the text (e.g. "41.40" -> 41.4,
# "33.33" -> 33.329999999999998...). Force each of those cells to Text format
# first so the literal digits are preserved exactly, then restore the original
# (unstyled) look once the values are in place.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D33", "D35", "D37", "D39", "D44", "D45", "D47", "D48", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin price / volume(1h) figures
$ws.Range("D2").Value = "60.822.60"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "2.911.14"
$ws.Range("D5").Value = "586.94"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "144.42"
$ws.Range("E6").Value = "  -6.24%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("D9").Value = "2.909.76"
$ws.Range("E9").Value = "  -3.94%  "
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("E11").Value = "  -5.38%  "
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("D14").Value = "33.33"
$ws.Range("E14").Value = "  -6.88%  "
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "3.395.40"
$ws.Range("E16").Value = "  -4.05%  "
$ws.Range("D17").Value = "60.780.51"
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -5.12%  "
$ws.Range("D19").Value = "2.912.82"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").Value = "427.99"
$ws.Range("E20").Value = "  -5.66%  "
$ws.Range("D21").Value = "13.53"
$ws.Range("E21").Value = "  -5.29%  "
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "7.09"
$ws.Range("D24").Value = "80.58"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").Value = "10.71"
$ws.Range("E26").Value = "  -4.49%  "
$ws.Range("D27").Value = "11.92"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").Value = "26.47"
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("D35").Value = "0.0₃0873"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  -2.76%  "
$ws.Range("D37").Value = "5.60"
$ws.Range("E37").Value = "  -5.62%  "
$ws.Range("E38").Value = "  -6.58%  "
$ws.Range("D39").Value = "49.50"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("E41").Value = "  -5.79%  "
$ws.Range("E42").Value = "  -5.95%  "
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").Value = "41.40"
$ws.Range("E44").Value = "  -5.25%  "
$ws.Range("D45").Value = "377.24"
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "2.682.25"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "132.29"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D50").Value = "24.42"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -2.63%  "

# Restore the default (no explicit style) look for the price cells
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
